# Add an "email" column to the voters sheet, right after "fullname" and
# before "voteid"/"voted" (the existing columns shift one place to the
# right), matching the commit "Add email address to voter record."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("voters")

# Insert a new column D; existing D (voteid) and E (voted) shift to E/F.
$ws.Columns.Item(4).Insert()

# Give the new column its header text.
$ws.Cells.Item(1, 4).Value = "email"

# Reflect that "voters" is now the sheet the user is working on.
$ws.Activate()
$null = $ws.Range("D2").Select()
